$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "w02"
$ws.Range("F3").Value = "w02"
$ws.Range("G3").Value = "w02"
$ws.Range("H3").Value = "w02"
$ws.Range("I3").Value = "w02"
$ws.Range("P6").Value = "w02"
$ws.Range("Q6").Value = "w02"
$ws.Range("R6").Value = "w02"
$ws.Range("S6").Value = "w02"
$ws.Range("T6").Value = "w02"
$ws.Range("AM8").Value = "w02"
$ws.Range("AN8").Value = "w02"
$ws.Range("AO8").Value = "w02"
$ws.Range("AQ8").Value = "w02"
$ws.Range("AR8").Value = "w02"
$ws.Range("AS8").Value = "w02"
$ws.Range("D9").Value = "w02"
$ws.Range("E9").Value = "w02"
$ws.Range("F9").Value = "w02"
$ws.Range("H9").Value = "w02"
$ws.Range("I9").Value = "w02"
$ws.Range("J9").Value = "w02"
$ws.Range("K9").Value = "w02"
$ws.Range("L9").Value = "w02"
$ws.Range("M9").Value = "w02"
$ws.Range("N9").Value = "w02"
$ws.Range("O9").Value = "w02"
$ws.Range("D14").Value = "w02"
$ws.Range("E14").Value = "w02"
$ws.Range("F14").Value = "w02"
$ws.Range("H14").Value = "w02"
$ws.Range("I14").Value = "w02"
$ws.Range("P18").Value = "w02"
$ws.Range("Q18").Value = "w02"
$ws.Range("S18").Value = "w02"
$ws.Range("T18").Value = "w02"
$ws.Range("Q23").Value = "d12"
$ws.Range("N24").Value = "w02"
$ws.Range("O24").Value = "w02"
$ws.Range("P24").Value = "w02"
$ws.Range("Q24").Value = "w02"
$ws.Range("S24").Value = "w02"
$ws.Range("T24").Value = "w02"
$ws.Range("AE24").Value = "d07"
$ws.Range("E25").Value = "w02"
$ws.Range("F25").Value = "w02"
$ws.Range("G25").Value = "w02"
$ws.Range("H25").Value = "w02"
$ws.Range("I25").Value = "w02"
$ws.Range("J25").Value = "w02"
$ws.Range("K25").Value = "w02"
$ws.Range("L25").Value = "w02"
$ws.Range("M25").Value = "w02"
$ws.Range("U25").Value = "w02"
$ws.Range("V25").Value = "w02"
$ws.Range("W25").Value = "w02"
$ws.Range("X25").Value = "w02"
$ws.Range("Y25").Value = "w02"
$ws.Range("Z25").Value = "w02"
$ws.Range("AA25").Value = "w02"
$ws.Range("AB25").Value = "w02"
$ws.Range("AC25").Value = "w02"
$ws.Range("AD25").Value = "w02"
$ws.Range("O28").Value = "d02"
$ws.Range("Q28").Value = "d03"
$ws.Range("AE28").Value = "d13"
$ws.Range("AN28").Value = "w02"
$ws.Range("AO28").Value = "w02"
$ws.Range("AQ28").Value = "w02"
$ws.Range("AR28").Value = "w02"
$ws.Range("AC29").Value = "w02"
$ws.Range("AE29").Value = "w02"
$ws.Range("AF29").Value = "w02"
$ws.Range("AG29").Value = "w02"
$ws.Range("O30").Value = "w02"
$ws.Range("Q30").Value = "w02"
$ws.Range("R30").Value = "w02"
$ws.Range("S31").Value = "w02"
$ws.Range("T31").Value = "w02"
$ws.Range("U31").Value = "w02"
$ws.Range("V31").Value = "w02"
$ws.Range("W31").Value = "w02"
$ws.Range("X31").Value = "w02"
$ws.Range("Y31").Value = "w02"
$ws.Range("Q35").Value = "d03"
$ws.Range("Q36").Value = "d03"
$ws.Range("Q37").Value = "d03"
$ws.Range("Q38").Value = "d13"
$ws.Range("AQ38").Value = "w02"
$ws.Range("AR38").Value = "w02"
$ws.Range("AS38").Value = "w02"
$ws.Range("AT38").Value = "w02"
$ws.Range("AU38").Value = "w02"
$ws.Range("AV38").Value = "w02"
$ws.Range("AW38").Value = "w02"
$ws.Range("R39").Value = "w02"
$ws.Range("S39").Value = "w02"
$ws.Range("T39").Value = "w02"
$ws.Range("U39").Value = "w02"
$ws.Range("V39").Value = "w02"
$ws.Range("W39").Value = "w02"
$ws.Range("X39").Value = "w02"
$ws.Range("D43").Value = "w02"
$ws.Range("E43").Value = "w02"
$ws.Range("F43").Value = "w02"
$ws.Range("H43").Value = "w02"
$ws.Range("I43").Value = "w02"
$ws.Range("J43").Value = "w02"
$ws.Range("K43").Value = "w02"
$ws.Range("L43").Value = "w02"
$ws.Range("S44").Value = "w02"
$ws.Range("T44").Value = "w02"
$ws.Range("U44").Value = "w02"
$ws.Range("V44").Value = "w02"
$ws.Range("W44").Value = "w02"
$ws.Range("M45").Value = "w02"
$ws.Range("N45").Value = "w02"
$ws.Range("O45").Value = "w02"
$ws.Range("P45").Value = "w02"
$ws.Range("Q45").Value = "w02"
$ws.Range("R45").Value = "w02"
$ws.Range("Y45").Value = "w02"
$ws.Range("Z45").Value = "w02"
$ws.Range("AB48").Value = "w02"
$ws.Range("AC48").Value = "w02"
$ws.Range("AD48").Value = "w02"
$ws.Range("AE48").Value = "w02"
$ws.Range("AF48").Value = "w02"
$ws.Range("AG48").Value = "w02"
$ws.Range("AH49").Value = "w02"
$ws.Range("AI49").Value = "w02"
$ws.Range("AJ49").Value = "w02"
$ws.Range("AK49").Value = "w02"
$ws.Range("AL49").Value = "w02"
$ws.Range("AM49").Value = "w02"
$ws.Range("AN49").Value = "w02"
$ws.Range("AO49").Value = "w02"
$ws.Range("J52").Value = "w02"
$ws.Range("K52").Value = "w02"
$ws.Range("L52").Value = "w02"
$ws.Range("M52").Value = "w02"
$ws.Range("N52").Value = "w02"
$ws.Range("O52").Value = "w02"
$ws.Range("P52").Value = "w02"
$ws.Range("Q52").Value = "w02"
$ws.Range("R52").Value = "w02"
$ws.Range("S52").Value = "w02"
$ws.Range("T52").Value = "w02"
$ws.Range("U52").Value = "w02"
$ws.Range("V52").Value = "w02"
$ws.Range("W52").Value = "w02"
$ws.Range("X52").Value = "w02"
$ws.Range("Y52").Value = "w02"
$ws.Range("Z52").Value = "w02"
$ws.Range("AH53").Value = "w02"
$ws.Range("AI53").Value = "w02"
$ws.Range("AJ53").Value = "w02"
$ws.Range("AK53").Value = "w02"
$ws.Range("AL53").Value = "w02"
$ws.Range("AM53").Value = "w02"
$ws.Range("AN53").Value = "w02"
$ws.Range("AO53").Value = "w02"
$ws.Range("AP53").Value = "w02"
$ws.Range("AQ53").Value = "w02"
$ws.Range("AR53").Value = "w02"
$ws.Range("Q56").Value = "w02"
$ws.Range("R56").Value = "w02"
$ws.Range("S56").Value = "w02"
$ws.Range("T56").Value = "w02"
$ws.Range("U56").Value = "w02"
$ws.Range("V56").Value = "w02"
$ws.Range("W56").Value = "w02"
$ws.Range("X56").Value = "w02"
$ws.Range("Y56").Value = "w02"
$ws.Range("Z56").Value = "w02"
$ws.Range("AK57").Value = "w02"
$ws.Range("AL57").Value = "w02"
$ws.Range("AM57").Value = "w02"
$ws.Range("AN57").Value = "w02"
$ws.Range("AO58").Value = "w02"
$ws.Range("AP58").Value = "w02"
$ws.Range("AQ58").Value = "w02"
$ws.Range("R62").Value = "w02"
$ws.Range("S62").Value = "w02"
$ws.Range("T62").Value = "w02"
$ws.Range("U62").Value = "w02"
$ws.Range("V62").Value = "w02"
$ws.Range("W62").Value = "w02"
$ws.Range("X62").Value = "w02"
$ws.Range("Y62").Value = "w02"
$ws.Range("Z62").Value = "w02"
$ws.Range("AA62").Value = "w02"
$ws.Range("AG63").Value = "w02"
$ws.Range("AH63").Value = "w02"
$ws.Range("AI63").Value = "w02"
$ws.Range("AJ63").Value = "w02"
$ws.Range("AK63").Value = "w02"
$ws.Range("AL63").Value = "w02"
$ws.Range("AM63").Value = "w02"
$ws.Range("AN63").Value = "w02"
$ws.Range("AO63").Value = "w02"
$ws.Range("AP63").Value = "w02"
$ws.Range("AQ63").Value = "w02"
$ws.Range("AR63").Value = "w02"
$ws.Range("AS63").Value = "w02"
$ws.Range("AT63").Value = "w02"
$ws.Range("AU63").Value = "w02"
$ws.Range("AV63").Value = "w02"

$excel.ActiveWindow.Zoom = 55
$ws.Range("A1:AZ65").Select()
